$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.390.59"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.645.94"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'598.39"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "'154.65"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("D9").Value = "2.645.15"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "'0.144"
$ws.Range("E10").Value = "  +7.78%  "
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "'5.26"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "3.122.08"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "68.281.36"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "2.666.63"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").Value = "'11.46"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").Value = "'365.14"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("E22").Value = "  +3.96%  "
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").Value = "'74.15"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'9.87"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("D29").Value = "2.774.62"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").Value = "'575.61"
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("D32").Value = "'8.16"
$ws.Range("E32").Value = "  +4.76%  "
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("E35").Value = "  +4.80%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "'1.59"
$ws.Range("E37").Value = "  +5.98%  "
$ws.Range("D38").Value = "'159.61"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("E42").Value = "  +2.75%  "
$ws.Range("D43").Value = "0.0₆0341"
$ws.Range("E43").Value = "  +11.74%  "
$ws.Range("D44").Value = "'2.67"
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("E45").Value = "  +3.71%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'40.34"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'157.35"
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("E50").Value = "  +2.27%  "
$ws.Range("D51").Value = "'21.97"
$ws.Range("E51").Value = "  +2.75%  "
